# Correction of perc value
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Activate()

# Fix the 2030 B-column percentage (was fat-fingered as 0.01, should be 0.1)
$ws.Range("B9").Value = 0.1

# Fix the 2035 D-column percentage (was 1.6E-2, should be 0.16)
$ws.Range("D10").Value = 0.16

# Leave the selection where the user last clicked while reviewing the fix
$ws.Range("E10").Select()
